$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in column A (rows 2-10 all share the same date string)
$ws.Range("A2:A10").Value = "16.05.2025"

# Update the menu item descriptions in column B
$ws.Range("B2").Value = "Ryba miruna, domowe frytki, surówka + zupa kapuśniak ze słodkiej kapusty lub rosół"
$ws.Range("B3").Value = "Zrazy wieprzowe z mięsa mielonego, ziemniaki, surówka + zupa kapuśniak ze słodkiej kapusty lub rosół"
$ws.Range("B4").Value = "Ryba miruna, domowe frytki, surówka "
$ws.Range("B5").Value = "Zrazy wieprzowe z mięsa mielonego, ziemniaki, surówka "
$ws.Range("B6").Value = "Zupa kapuśniak ze słodkiej kapusty lub rosół"

# Update prices in column C
$ws.Range("C2").Value = 39
$ws.Range("C3").Value = 32
$ws.Range("C4").Value = 36
$ws.Range("C5").Value = 29

# Update the active cell selection to B6
$ws.Range("B6").Select()
